$wb = $excel.ActiveWorkbook

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3103.6667
$ws.Range("I62").Value = 1155.5
$ws.Range("K62").Value = 1155.5
$ws.Range("M62").Value = -531.5

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3103.6667
$ws.Range("I65").Value = 1155.5
$ws.Range("K65").Value = 5777.5
$ws.Range("M65").Value = -2657.5

# Sheet ALC, row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496

# Sheet ALC, row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 798.4915
$ws.Range("J132").Value = 1874.25
$ws.Range("L132").Value = 5622.75
$ws.Range("N132").Value = -10682.75

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3352.8223
$ws.Range("I138").Value = 1228.3
$ws.Range("J138").Value = 3959.8286
$ws.Range("K138").Value = 3684.9
$ws.Range("L138").Value = 11879.4858
$ws.Range("M138").Value = 1455.1
$ws.Range("N138").Value = -22159.4858

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1766468.6
$ws.Range("I32").Value = 8750.729499999999
$ws.Range("J32").Value = 16707071
$ws.Range("K32").Value = 8750.729499999999
$ws.Range("L32").Value = 16707071
$ws.Range("M32").Value = -8463.729499999999
$ws.Range("N32").Value = -16707645

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4113.6553
$ws.Range("I61").Value = 2413.4666
$ws.Range("J61").Value = 5935.2856
$ws.Range("K61").Value = 2413.4666
$ws.Range("L61").Value = 5935.2856
$ws.Range("M61").Value = -2201.4666
$ws.Range("N61").Value = -6359.2856

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1768.3658
$ws.Range("I74").Value = 1663.5264
$ws.Range("K74").Value = 1663.5264
$ws.Range("M74").Value = -789.5264

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1768.3658
$ws.Range("I77").Value = 1663.5264
$ws.Range("K77").Value = 8317.632
$ws.Range("M77").Value = -3949.632

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4113.6553
$ws.Range("I136").Value = 2413.4666
$ws.Range("J136").Value = 5935.2856
$ws.Range("K136").Value = 7240.399800000001
$ws.Range("L136").Value = 17805.8568
$ws.Range("M136").Value = -4690.399800000001
$ws.Range("N136").Value = -22905.8568

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4196.654
$ws.Range("I31").Value = 3147.9375
$ws.Range("K31").Value = 3147.9375
$ws.Range("M31").Value = -2852.9375

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4196.654
$ws.Range("I34").Value = 3147.9375
$ws.Range("K34").Value = 3147.9375
$ws.Range("M34").Value = -2945.9375

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2433.1553
$ws.Range("I134").Value = 1285.9166
$ws.Range("J134").Value = 7939.9
$ws.Range("K134").Value = 3857.7498
$ws.Range("L134").Value = 23819.7
$ws.Range("M134").Value = -1322.7498
$ws.Range("N134").Value = -28889.7

# Sheet CUL, row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 448.79166
$ws.Range("I7").Value = 293.41177
$ws.Range("J7").Value = 826.1429000000001
$ws.Range("K7").Value = 880.23531
$ws.Range("L7").Value = 2478.4287
$ws.Range("M7").Value = -768.23531
$ws.Range("N7").Value = -2702.4287

# Sheet CUL, row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4499.4546
$ws.Range("I64").Value = 3166.6667
$ws.Range("J64").Value = 4999.25
$ws.Range("K64").Value = 9500.000100000001
$ws.Range("L64").Value = 14997.75
$ws.Range("M64").Value = -9230.000100000001
$ws.Range("N64").Value = -15537.75

# Sheet CUL, row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 4499.4546
$ws.Range("I67").Value = 3166.6667
$ws.Range("J67").Value = 4999.25
$ws.Range("K67").Value = 9500.000100000001
$ws.Range("L67").Value = 14997.75
$ws.Range("M67").Value = -8564.000100000001
$ws.Range("N67").Value = -16869.75

# Sheet CUL, row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1987.25
$ws.Range("J75").Value = 1987.25
$ws.Range("L75").Value = 5961.75
$ws.Range("N75").Value = -7957.75

# Sheet CUL, row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1987.25
$ws.Range("J78").Value = 1987.25
$ws.Range("L78").Value = 17885.25
$ws.Range("N78").Value = -27869.25

# Sheet CUL, row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 3775
$ws.Range("I97").Value = 2800
$ws.Range("J97").Value = 4506.25
$ws.Range("K97").Value = 8400
$ws.Range("L97").Value = 13518.75
$ws.Range("M97").Value = -7904
$ws.Range("N97").Value = -14510.75

# Sheet CUL, row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 5000
$ws.Range("J124").Value = 5000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -24820

# Sheet CUL, row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 2497.5
$ws.Range("I126").Value = 2497.5
$ws.Range("K126").Value = 7492.5
$ws.Range("M126").Value = -2552.5

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3255.3438
$ws.Range("J131").Value = 3654.15
$ws.Range("L131").Value = 10962.45
$ws.Range("N131").Value = -21042.45

# Sheet GSM, row 48
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 19999
$ws.Range("I48").Value = 19999
$ws.Range("K48").Value = 19999
$ws.Range("M48").Value = -19514

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4473.5957
$ws.Range("I132").Value = 3779.4
$ws.Range("K132").Value = 11338.2
$ws.Range("M132").Value = -8808.200000000001

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4332.619
$ws.Range("I7").Value = 2523.5
$ws.Range("J7").Value = 5977.273
$ws.Range("K7").Value = 2523.5
$ws.Range("L7").Value = 5977.273
$ws.Range("M7").Value = -2411.5
$ws.Range("N7").Value = -6201.273

# Sheet LTW, row 95
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4332.619
$ws.Range("I126").Value = 2523.5
$ws.Range("J126").Value = 5977.273
$ws.Range("K126").Value = 7570.5
$ws.Range("L126").Value = 17931.819
$ws.Range("M126").Value = -5100.5
$ws.Range("N126").Value = -22871.819

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5307.7085
$ws.Range("I136").Value = 3399.375
$ws.Range("J136").Value = 9124.375
$ws.Range("K136").Value = 10198.125
$ws.Range("L136").Value = 27373.125
$ws.Range("M136").Value = -7648.125
$ws.Range("N136").Value = -32473.125

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 37039604
$ws.Range("I62").Value = 3849.5
$ws.Range("K62").Value = 3849.5
$ws.Range("M62").Value = -3225.5

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 37039604
$ws.Range("I65").Value = 3849.5
$ws.Range("K65").Value = 19247.5
$ws.Range("M65").Value = -16127.5

# Sheet WVR, row 97
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 59786
$ws.Range("J97").Value = 59786
$ws.Range("L97").Value = 59786
$ws.Range("N97").Value = -61768

# Sheet WVR, row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 58928.5
$ws.Range("J135").Value = 58928.5
$ws.Range("L135").Value = 58928.5
$ws.Range("N135").Value = -69068.5

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11114266
$ws.Range("I136").Value = 14287918
$ws.Range("J136").Value = 6484.9
$ws.Range("K136").Value = 42863754
$ws.Range("L136").Value = 19454.7
$ws.Range("M136").Value = -42861204
$ws.Range("N136").Value = -24554.7
